$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "Aluslevypari  M8, NL8SP_ExcelDescCol_IMP"
$ws.Range("E9").Select()
